$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SPAIN_Components_pop_change")

# For each "year block" in the components-of-population-change table, add the
# residual net migration figure for the block's total/summary row (total
# population change minus the natural increase for that block), mirroring
# the formula pattern already used on the existing summary rows such as
# C16, C25, C35 and C143: C(row-4) - C(row-1)
$rows = @(44, 53, 62, 71, 80, 89, 98, 107, 116, 125, 134)

foreach ($r in $rows) {
    $dst = $ws.Range("C$r")

    # Pick up the same cell style/number format used by the other
    # "Net migration" total rows in this table (e.g. C35) before writing
    # the new formula, so the new cells render identically to their peers.
    $ws.Range("C35").Copy()
    $dst.PasteSpecial(-4122)  # xlPasteFormats

    $dst.Formula = "=C$($r - 4)-C$($r - 1)"
}

$excel.CutCopyMode = 0

# Put the cursor on the last cell touched while entering the formulas above.
$ws.Range("C143").Select()
